$wb = $excel.ActiveWorkbook

# Rename the "October 2022" sheet to "November 2022"
$ws = $wb.Worksheets.Item("October 2022")
$ws.Name = "November 2022"

# Update data rows: keep rows 1 (header) and 2, overwrite with new calibration data, delete rows 4-6
$ws.Range("A2").Value = "LBA51"
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = "10/13"

$ws.Range("A3").Value = "LBC62"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = "10/13"

# Remove old rows 4-6 (shift up / clear)
$ws.Range("A4:F6").Delete()

# Update selection to K24 like in the diff
$ws.Range("K24").Select()
